$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17
$ws.Range("H17").Value = 307.84506
$ws.Range("J17").Value = 307.84506
$ws.Range("L17").Value = 923.53518
$ws.Range("N17").Value = -1259.53518

# ALC row 127
$ws.Range("H127").Value = 21741302
$ws.Range("I127").Value = 424.2857
$ws.Range("J127").Value = 25643510
$ws.Range("K127").Value = 1272.8571
$ws.Range("L127").Value = 76930530
$ws.Range("M127").Value = 3687.1429
$ws.Range("N127").Value = -76940450

# ALC row 129
$ws.Range("H129").Value = 290182.56
$ws.Range("J129").Value = 309490.34
$ws.Range("L129").Value = 928471.02
$ws.Range("N129").Value = -938471.02

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Range("H61").Value = 3267.7778
$ws.Range("J61").Value = 3999.6667
$ws.Range("L61").Value = 3999.6667
$ws.Range("N61").Value = -4423.6667

# ARM row 132
$ws.Range("H132").Value = 1708.5834
$ws.Range("I132").Value = 1336.1072
$ws.Range("J132").Value = 3012.25
$ws.Range("K132").Value = 4008.3216
$ws.Range("L132").Value = 9036.75
$ws.Range("M132").Value = -1478.3216
$ws.Range("N132").Value = -14096.75

# ARM row 136
$ws.Range("H136").Value = 3267.7778
$ws.Range("J136").Value = 3999.6667
$ws.Range("L136").Value = 11999.0001
$ws.Range("N136").Value = -17099.0001

$ws = $wb.Worksheets.Item("BSM")
# BSM row 76
$ws.Range("H76").Value = 28550
$ws.Range("J76").Value = 28550
$ws.Range("L76").Value = 28550
$ws.Range("N76").Value = -29180

# BSM row 79
$ws.Range("H79").Value = 28550
$ws.Range("J79").Value = 28550
$ws.Range("L79").Value = 28550
$ws.Range("N79").Value = -30734

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 31835.941
$ws.Range("I31").Value = 53989.21
$ws.Range("J31").Value = 3775.1333
$ws.Range("K31").Value = 53989.21
$ws.Range("L31").Value = 3775.1333
$ws.Range("M31").Value = -53694.21
$ws.Range("N31").Value = -4365.1333

# CRP row 34
$ws.Range("H34").Value = 31835.941
$ws.Range("I34").Value = 53989.21
$ws.Range("J34").Value = 3775.1333
$ws.Range("K34").Value = 53989.21
$ws.Range("L34").Value = 3775.1333
$ws.Range("M34").Value = -53787.21
$ws.Range("N34").Value = -4179.1333

# CRP row 69
$ws.Range("H69").Value = 17500
$ws.Range("I69").Value = 17500
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 17500
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -16751
$ws.Range("N69").ClearContents()

# CRP row 72
$ws.Range("H72").Value = 17500
$ws.Range("I72").Value = 17500
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 52500
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -48756
$ws.Range("N72").ClearContents()

# CRP row 92
$ws.Range("H92").Value = 36196
$ws.Range("J92").Value = 36196
$ws.Range("L92").Value = 36196
$ws.Range("N92").Value = -41188

# CRP row 134
$ws.Range("H134").Value = 1394.4445
$ws.Range("I134").Value = 1346
$ws.Range("K134").Value = 4038
$ws.Range("M134").Value = -1503

$ws = $wb.Worksheets.Item("CUL")
# CUL row 37
$ws.Range("H37").Value = 579061.3
$ws.Range("J37").Value = 579061.3
$ws.Range("L37").Value = 1737183.9
$ws.Range("N37").Value = -1737407.9

# CUL row 131
$ws.Range("H131").Value = 845.4
$ws.Range("I131").Value = 513.1667
$ws.Range("J131").Value = 866.6064
$ws.Range("K131").Value = 1539.5001
$ws.Range("L131").Value = 2599.8192
$ws.Range("M131").Value = 3500.4999
$ws.Range("N131").Value = -12679.8192

# CUL row 134
$ws.Range("H134").Value = 5243.909
$ws.Range("I134").Value = 3285
$ws.Range("J134").Value = 7594.6
$ws.Range("K134").Value = 9855
$ws.Range("L134").Value = 22783.8
$ws.Range("M134").Value = -4785
$ws.Range("N134").Value = -32923.8

# CUL row 136
$ws.Range("H136").Value = 2650
$ws.Range("I136").Value = 1600
$ws.Range("K136").Value = 4800
$ws.Range("M136").Value = 300

# CUL row 138
$ws.Range("H138").Value = 1933.25
$ws.Range("I138").Value = 1455.4445
$ws.Range("J138").Value = 3366.6667
$ws.Range("K138").Value = 4366.333500000001
$ws.Range("L138").Value = 10100.0001
$ws.Range("M138").Value = 773.6664999999994
$ws.Range("N138").Value = -20380.0001

# CUL row 139
$ws.Range("H139").Value = 2406.923
$ws.Range("I139").Value = 1589.1666
$ws.Range("J139").Value = 3107.8572
$ws.Range("K139").Value = 4767.4998
$ws.Range("L139").Value = 9323.571599999999
$ws.Range("M139").Value = 372.5002000000004
$ws.Range("N139").Value = -19603.5716

# CUL row 140
$ws.Range("H140").Value = 1559.3103
$ws.Range("I140").Value = 1024.7059
$ws.Range("J140").Value = 2316.6667
$ws.Range("K140").Value = 3074.1177
$ws.Range("L140").Value = 6950.000100000001
$ws.Range("M140").Value = 2105.8823
$ws.Range("N140").Value = -17310.0001

# CUL row 141
$ws.Range("H141").Value = 3210.6667
$ws.Range("I141").Value = 2892.8
$ws.Range("J141").Value = 4800
$ws.Range("K141").Value = 8678.400000000001
$ws.Range("L141").Value = 14400
$ws.Range("M141").Value = -3498.400000000001
$ws.Range("N141").Value = -24760

$ws = $wb.Worksheets.Item("GSM")
# GSM row 17
$ws.Range("H17").Value = 11950
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 11950
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 11950
$ws.Range("N17").Value = -12286
$ws.Range("M17").ClearContents()

# GSM row 70
$ws.Range("H70").Value = 61688.855
$ws.Range("I70").Value = 104206.9
$ws.Range("J70").Value = 4998.1333
$ws.Range("K70").Value = 104206.9
$ws.Range("L70").Value = 4998.1333
$ws.Range("M70").Value = -103936.9
$ws.Range("N70").Value = -5538.1333

# GSM row 73
$ws.Range("H73").Value = 61688.855
$ws.Range("I73").Value = 104206.9
$ws.Range("J73").Value = 4998.1333
$ws.Range("K73").Value = 104206.9
$ws.Range("L73").Value = 4998.1333
$ws.Range("M73").Value = -103270.9
$ws.Range("N73").Value = -6870.1333

# GSM row 122
$ws.Range("H122").Value = 1901.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1901.25
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 5703.75
$ws.Range("N122").Value = -10603.75
$ws.Range("M122").ClearContents()

# GSM row 126
$ws.Range("H126").Value = 2104.2917
$ws.Range("I126").Value = 1742.6666
$ws.Range("J126").Value = 2321.2666
$ws.Range("K126").Value = 5227.9998
$ws.Range("L126").Value = 6963.7998
$ws.Range("M126").Value = -2757.9998
$ws.Range("N126").Value = -11903.7998

$ws = $wb.Worksheets.Item("LTW")
# LTW row 122
$ws.Range("H122").Value = 2266.3333
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2266.3333
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6798.999899999999
$ws.Range("N122").Value = -11698.9999
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# WVR row 69
$ws.Range("H69").Value = 21940
$ws.Range("J69").Value = 21940
$ws.Range("L69").Value = 21940
$ws.Range("N69").Value = -23438

# WVR row 72
$ws.Range("H72").Value = 21940
$ws.Range("J72").Value = 21940
$ws.Range("L72").Value = 65820
$ws.Range("N72").Value = -73308
